$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the target formatting (custom 10-decimal number format with thousands
# separator, plus a thin box border) on a scratch cell off to the side, then
# copy/paste the *formats only* onto the data range in a single atomic
# operation. Doing it this way (rather than touching NumberFormat/Borders
# directly on the destination range) avoids leaving an extra, unused
# intermediate cell style behind in the workbook's style table.
$scratch = $ws.Range("D1")
$scratch.NumberFormat = "#,##0.0000000000"
$scratch.Borders.LineStyle = 1
$scratch.Copy()

$dataRange = $ws.Range("A1:B7")
$dataRange.PasteSpecial(-4122)  # xlPasteFormats

$scratch.Clear()
$excel.CutCopyMode = 0

# Refresh the data values (new sample, regenerated after the decimal point fix)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 0

$ws.Range("A2").Value = -1.917462659138496
$ws.Range("B2").Value = -0.5686272512019883

$ws.Range("A3").Value = -2.011525851630627
$ws.Range("B3").Value = -3.735473697833986

$ws.Range("A4").Value = 0.8528022432798865
$ws.Range("B4").Value = -2.876235097111638

$ws.Range("A5").Value = 1.196765323619466
$ws.Range("B5").Value = -4.0675979103375

$ws.Range("A6").Value = -1.24293612056506
$ws.Range("B6").Value = 0.6746182617891435

$ws.Range("A7").Value = 1.242341834377565
$ws.Range("B7").Value = -0.6757120433291937
